$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 52: continuation of Hugo's test-sound catalogue entries
$ws.Range("A52").Value = "Hugo"
$ws.Range("B52").Value = "Ni"
$ws.Range("D52").Value = "underground                                                                      /  city"

# Move the viewport/selection to show the newly added row, like the author did
$ws.Range("J53").Select() | Out-Null
